# Redefine definition of financial cost as effective cost/loan ratio
# Update the hidden helper columns (I:L) on each results sheet with the
# recomputed cost/loan-ratio figures. The visible B:E columns are driven by
# ROUND() formulas over I:L and are intentionally left to recalc on next open
# (Calculation kept Manual here) rather than recalculated inline now.
$excel.Calculation = -4135  # xlCalculationManual: keep B:E caches as-is
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("I4").Value = 0.29876092299908968
$ws.Range("J4").Value = 0.29995205830602678
$ws.Range("K4").Value = 0.31539230728926859
$ws.Range("L4").Value = 0.32024212869037177
$ws.Range("I5").Value = 0.14922114386074789
$ws.Range("J5").Value = 0.1509066226977657
$ws.Range("K5").Value = 0.15150041626271421
$ws.Range("L5").Value = 0.15565407794454819
$ws.Range("I6").Value = 0.77375795412879844
$ws.Range("J6").Value = 0.76484076006017843
$ws.Range("K6").Value = 0.79757961988638904
$ws.Range("L6").Value = 0.77057324151514384
$ws.Range("I7").Value = 0.039567531636237198
$ws.Range("J7").Value = 0.040175793327711498
$ws.Range("K7").Value = 0.037045175218860499
$ws.Range("L7").Value = 0.038629866030232697
$ws.Range("I8").Value = 0.75927400673325851
$ws.Range("J8").Value = 0.75221204142830833
$ws.Range("K8").Value = 0.87379146035638344
$ws.Range("L8").Value = 0.83834755707842912
$ws.Range("I9").Value = 0.015715355688691601
$ws.Range("J9").Value = 0.015789085696552501
$ws.Range("K9").Value = 0.013753055858661101
$ws.Range("L9").Value = 0.0134523148749111
$ws.Range("I10").Value = 0.77294685990338163
$ws.Range("J10").Value = 0.78260869565217395
$ws.Range("K10").Value = 0.75362318840579712
$ws.Range("L10").Value = 0.77777777777777779
$ws.Range("I11").Value = 0.40155597847888153
$ws.Range("J11").Value = 0.39462577544128402
$ws.Range("K11").Value = 0.36875261807523541
$ws.Range("L11").Value = 0.38525974238462968
$ws.Range("I12").Value = 0.43156871865542351
$ws.Range("J12").Value = 0.4216236480314085
$ws.Range("K12").Value = 0.4496237378390267
$ws.Range("L12").Value = 0.41262911550148101
$ws.Range("I13").Value = 0.18941553929808749
$ws.Range("J13").Value = 0.18025984648686291
$ws.Range("K13").Value = 0.17703429932806869
$ws.Range("L13").Value = 0.15447118576924521
$ws.Range("I14").Value = 0.23542725189098559
$ws.Range("J14").Value = 0.23046351785227151
$ws.Range("K14").Value = 0.23221790441424381
$ws.Range("L14").Value = 0.2321760727155612

$ws = $wb.Worksheets.Item(2)
$ws.Range("I4").Value = 0.2227522761263121
$ws.Range("J4").Value = 0.21598377519419321
$ws.Range("K4").Value = 0.23237049240652499
$ws.Range("L4").Value = 0.21188329852202931
$ws.Range("I5").Value = 0.12074442167539851
$ws.Range("J5").Value = 0.11169839120125991
$ws.Range("K5").Value = 0.1072400325053447
$ws.Range("L5").Value = 0.11045445622624581
$ws.Range("I6").Value = 0.76388889542884253
$ws.Range("J6").Value = 0.76875000440826036
$ws.Range("K6").Value = 0.78055555888244665
$ws.Range("L6").Value = 0.78680556117970901
$ws.Range("I7").Value = 0.053243170521964998
$ws.Range("J7").Value = 0.0541617061220648
$ws.Range("K7").Value = 0.051753578782367297
$ws.Range("L7").Value = 0.051569759765743399
$ws.Range("I8").Value = 0.83742630419953501
$ws.Range("J8").Value = 0.82642409117787263
$ws.Range("K8").Value = 0.9296684260848973
$ws.Range("L8").Value = 0.94874394949001783
$ws.Range("I9").Value = 0.0200607760744052
$ws.Range("J9").Value = 0.020491937260801701
$ws.Range("K9").Value = 0.0136267352260047
$ws.Range("L9").Value = 0.0107578042101512
$ws.Range("I10").Value = 0.78378378378378377
$ws.Range("J10").Value = 0.81081081081081086
$ws.Range("K10").Value = 0.78378378378378377
$ws.Range("L10").Value = 0.81081081081081086
$ws.Range("I11").Value = 0.17018128710451469
$ws.Range("J11").Value = 0.25505639266809238
$ws.Range("K11").Value = 0.12459257757982591
$ws.Range("L11").Value = 0.21291561523412569
$ws.Range("I12").Value = 0.30655431494183077
$ws.Range("J12").Value = 0.32088049258398949
$ws.Range("K12").Value = 0.31842320254232231
$ws.Range("L12").Value = 0.31371607570461663
$ws.Range("I13").Value = -0.023812075455982601
$ws.Range("J13").Value = 0.0528907246059842
$ws.Range("K13").Value = 0.090693891048431896
$ws.Range("L13").Value = 0.063438256581624594
$ws.Range("I14").Value = 0.18243361448004911
$ws.Range("J14").Value = 0.1667835349736207
$ws.Range("K14").Value = 0.15269345887475189
$ws.Range("L14").Value = 0.1566764534857997

$ws = $wb.Worksheets.Item(3)
$ws.Range("I4").Value = 0.36381110744402212
$ws.Range("J4").Value = 0.37656269091716471
$ws.Range("K4").Value = 0.37071027848473248
$ws.Range("L4").Value = 0.39295454212325692
$ws.Range("I5").Value = 0.196185051126683
$ws.Range("J5").Value = 0.20775241040168921
$ws.Range("K5").Value = 0.18006483037421039
$ws.Range("L5").Value = 0.18908884170207571
$ws.Range("I6").Value = 0.71830237819479725
$ws.Range("J6").Value = 0.68647214099814191
$ws.Range("K6").Value = 0.76021219716217114
$ws.Range("L6").Value = 0.70265251479459934
$ws.Range("I7").Value = 0.055803131672913002
$ws.Range("J7").Value = 0.056787581879631802
$ws.Range("K7").Value = 0.053049616326186898
$ws.Range("L7").Value = 0.055729715648370599
$ws.Range("I8").Value = 0.81609761176367968
$ws.Range("J8").Value = 0.8012186582778521
$ws.Range("K8").Value = 0.89510896161019415
$ws.Range("L8").Value = 0.81547858781499705
$ws.Range("I9").Value = 0.0186405777666422
$ws.Range("J9").Value = 0.019234581565911599
$ws.Range("K9").Value = 0.014013606646773299
$ws.Range("L9").Value = 0.018837461106733198
$ws.Range("I10").Value = 0.69148936170212771
$ws.Range("J10").Value = 0.63829787234042556
$ws.Range("K10").Value = 0.7021276595744681
$ws.Range("L10").Value = 0.69148936170212771
$ws.Range("I11").Value = 0.2700373738300143
$ws.Range("J11").Value = 0.16834325623917759
$ws.Range("K11").Value = 0.31531786385356791
$ws.Range("L11").Value = 0.27003737383001419
$ws.Range("I12").Value = 0.34192434086511958
$ws.Range("J12").Value = 0.30924956867698988
$ws.Range("K12").Value = 0.40132745180685392
$ws.Range("L12").Value = 0.34824554293972498
$ws.Range("I13").Value = 0.090159301353070001
$ws.Range("J13").Value = 0.036513755049566402
$ws.Range("K13").Value = 0.16491949754305171
$ws.Range("L13").Value = 0.1230691489869782
$ws.Range("I14").Value = 0.34201023863638408
$ws.Range("J14").Value = 0.36477182341858427
$ws.Range("K14").Value = 0.34157019053050808
$ws.Range("L14").Value = 0.34778111982852861
